$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.750.03'
$ws.Range('E2').Value = '  -2.86%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.500.67'
$ws.Range('E3').Value = '  -0.25%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '557.30'
$ws.Range('E5').Value = '  +1.24%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '179.57'
$ws.Range('E6').Value = '  -5.54%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.641'
$ws.Range('E7').Value = '  +5.50%  '

$ws.Range('E8').Value = '  +0.14%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.631'
$ws.Range('E9').Value = '  -0.93%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.157'
$ws.Range('E10').Value = '  +4.01%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.96'
$ws.Range('E11').Value = '  -5.36%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000273'
$ws.Range('E12').Value = '  -0.06%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.16'
$ws.Range('E13').Value = '  -2.78%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.065.49'
$ws.Range('E14').Value = '  +0.01%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.510.00'
$ws.Range('E15').Value = '  +0.07%  '

$ws.Range('E16').Value = '  +0.21%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '18.39'
$ws.Range('E17').Value = '  +0.45%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.14'
$ws.Range('E18').Value = '  +2.06%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '65.805.60'
$ws.Range('E19').Value = '  -3.29%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.999'
$ws.Range('E20').Value = '  +0.01%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '414.34'
$ws.Range('E21').Value = '  +1.30%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.09'
$ws.Range('E22').Value = '  +3.78%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '85.43'
$ws.Range('E23').Value = '  +1.24%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.09'
$ws.Range('E24').Value = '  -2.57%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.74'
$ws.Range('E25').Value = '  +7.57%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.82'
$ws.Range('E26').Value = '  -6.43%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.84'
$ws.Range('E27').Value = '  -3.09%  '

$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.05'
$ws.Range('E28').Value = '  -1.34%  '

$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.07'
$ws.Range('E29').Value = '  +5.56%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '30.28'
$ws.Range('E30').Value = '  -0.67%  '

$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '619.61'
$ws.Range('E31').Value = '  -8.64%  '

$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.41'
$ws.Range('E32').Value = '  -6.06%  '

$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.62'
$ws.Range('E33').Value = '  -1.15%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.110'
$ws.Range('E34').Value = '  -1.04%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '59.21'
$ws.Range('E35').Value = '  -1.84%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.148'
$ws.Range('E36').Value = '  +11.36%  '

$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0₃0801'
$ws.Range('E38').Value = '  -2.90%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.402.56'
$ws.Range('E39').Value = '  +12.19%  '

$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '37.09'
$ws.Range('E40').Value = '  -5.10%  '

$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.380'
$ws.Range('E41').Value = '  -4.95%  '

$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.33'
$ws.Range('E42').Value = '  -4.03%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.90'
$ws.Range('E44').Value = '  -3.27%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.33'
$ws.Range('E45').Value = '  -1.49%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0418'
$ws.Range('E46').Value = '  -1.24%  '

$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.50'
$ws.Range('E47').Value = '  -7.37%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.74'
$ws.Range('E48').Value = '  -0.30%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.133'
$ws.Range('E49').Value = '  +1.73%  '

$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.45'
$ws.Range('E50').Value = '  -5.55%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '137.15'
$ws.Range('E51').Value = '  -1.40%  '
